{"js": "// US12 title update:\n// \"ID User Story : UC 12 Visualizar a descri\u00e7\u00e3o de cada item do plano\n//  facilitando o preenchimento do plano.\"\n// becomes\n// \"ID User Story : 12 Visualiza\u00e7\u00e3o de Descri\u00e7\u00e3o e Exemplos de Cada Item\n//  do Plano de Neg\u00f3cio.\"\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\n\n// 1) Drop the \"UC \" prefix in front of the case number (\"UC 12\" -> \"12\").\nlet results = titlePara.search(\"UC 12\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"12\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Replace the descriptive sentence with the new wording, keeping the\n//    existing bold/Arial formatting of the run it overwrites.\nresults = titlePara.search(\n  \"Visualizar a descri\u00e7\u00e3o de cada item do plano facilitando o preenchimento do plano.\",\n  { matchCase: true }\n);\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\n    \"Visualiza\u00e7\u00e3o de Descri\u00e7\u00e3o e Exemplos de Cada Item do Plano de Neg\u00f3cio.\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "# US12 title update:\n# \"ID User Story : UC 12 Visualizar a descri\u00e7\u00e3o de cada item do plano\n#  facilitando o preenchimento do plano.\"\n# becomes\n# \"ID User Story : 12 Visualiza\u00e7\u00e3o de Descri\u00e7\u00e3o e Exemplos de Cada Item\n#  do Plano de Neg\u00f3cio.\"\n$d = $word.ActiveDocument\n$wdReplaceAll = 2\n\n# 1) Drop the \"UC \" prefix in front of the case number (\"UC 12\" -> \"12\").\n$find1 = $d.Content.Find\n$find1.Text = \"UC 12\"\n$find1.Replacement.Text = \"12\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, $wdReplaceAll)\n\n# 2) Replace the descriptive sentence with the new wording, keeping the\n#    existing bold/Arial formatting of the range it overwrites.\n$find2 = $d.Content.Find\n$find2.Text = \"Visualizar a descri\u00e7\u00e3o de cada item do plano facilitando o preenchimento do plano.\"\n$find2.Replacement.Text = \"Visualiza\u00e7\u00e3o de Descri\u00e7\u00e3o e Exemplos de Cada Item do Plano de Neg\u00f3cio.\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, $wdReplaceAll)\n"}
